$wb = $excel.ActiveWorkbook

# Sheet 1 (index 1): "Exhibitions" sheet - update column F ("want to go" counts)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 35
$ws.Cells.Item(3, 6).Value = 66
$ws.Cells.Item(5, 6).Value = 8875
$ws.Cells.Item(6, 6).Value = 145
$ws.Cells.Item(7, 6).Value = 236
$ws.Cells.Item(8, 6).Value = 7101
$ws.Cells.Item(9, 6).Value = 184
$ws.Cells.Item(10, 6).Value = 5405
$ws.Cells.Item(11, 6).Value = 10
$ws.Cells.Item(12, 6).Value = 6139
$ws.Cells.Item(14, 6).Value = 398
$ws.Cells.Item(18, 6).Value = 328
$ws.Cells.Item(19, 6).Value = 271
$ws.Cells.Item(21, 6).Value = 205
$ws.Cells.Item(22, 6).Value = 161
$ws.Cells.Item(23, 6).Value = 104
$ws.Cells.Item(24, 6).Value = 10005
$ws.Cells.Item(25, 6).Value = 1903
$ws.Cells.Item(26, 6).Value = 1820
$ws.Cells.Item(27, 6).Value = 47
$ws.Cells.Item(28, 6).Value = 36
$ws.Cells.Item(29, 6).Value = 2061
$ws.Cells.Item(30, 6).Value = 80
$ws.Cells.Item(31, 6).Value = 86
$ws.Cells.Item(32, 6).Value = 162
$ws.Cells.Item(35, 6).Value = 305
$ws.Cells.Item(36, 6).Value = 1365
$ws.Cells.Item(37, 6).Value = 61
$ws.Cells.Item(38, 6).Value = 5161
$ws.Cells.Item(39, 6).Value = 393
$ws.Cells.Item(40, 6).Value = 1194
$ws.Cells.Item(41, 6).Value = 648
$ws.Cells.Item(42, 6).Value = 103
$ws.Cells.Item(45, 6).Value = 1074
$ws.Cells.Item(46, 6).Value = 975
$ws.Cells.Item(47, 6).Value = 1354
$ws.Cells.Item(48, 6).Value = 62
$ws.Cells.Item(49, 6).Value = 1093

# Sheet 2 (index 2): "Performances" sheet - update column F ("want to go" counts)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 37
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(4, 6).Value = 15
$ws.Cells.Item(7, 6).Value = 31
$ws.Cells.Item(9, 6).Value = 41
$ws.Cells.Item(10, 6).Value = 190
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(19, 6).Value = 7

# Sheet 4 (index 4): "All Types" sheet - update column F ("want to go" counts)
# (Sheet 3, "Local Life", has no changes in this update.)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 35
$ws.Cells.Item(3, 6).Value = 66
$ws.Cells.Item(4, 6).Value = 26
$ws.Cells.Item(5, 6).Value = 15
$ws.Cells.Item(6, 6).Value = 8875
$ws.Cells.Item(7, 6).Value = 145
$ws.Cells.Item(8, 6).Value = 236
$ws.Cells.Item(9, 6).Value = 7101
$ws.Cells.Item(10, 6).Value = 184
$ws.Cells.Item(11, 6).Value = 5
$ws.Cells.Item(13, 6).Value = 5405
$ws.Cells.Item(14, 6).Value = 6139
$ws.Cells.Item(15, 6).Value = 6139
$ws.Cells.Item(17, 6).Value = 398
$ws.Cells.Item(19, 6).Value = 31
$ws.Cells.Item(20, 6).Value = 560
$ws.Cells.Item(22, 6).Value = 271
$ws.Cells.Item(24, 6).Value = 205
$ws.Cells.Item(26, 6).Value = 104
$ws.Cells.Item(27, 6).Value = 190
$ws.Cells.Item(28, 6).Value = 10005
$ws.Cells.Item(29, 6).Value = 1903
$ws.Cells.Item(30, 6).Value = 1820
$ws.Cells.Item(31, 6).Value = 47
$ws.Cells.Item(32, 6).Value = 2061
$ws.Cells.Item(33, 6).Value = 80
$ws.Cells.Item(35, 6).Value = 1031
$ws.Cells.Item(37, 6).Value = 2045
$ws.Cells.Item(38, 6).Value = 305
$ws.Cells.Item(39, 6).Value = 1365
$ws.Cells.Item(40, 6).Value = 5161
$ws.Cells.Item(42, 6).Value = 648
$ws.Cells.Item(43, 6).Value = 103
$ws.Cells.Item(46, 6).Value = 1074
$ws.Cells.Item(48, 6).Value = 1354
$ws.Cells.Item(49, 6).Value = 62
$ws.Cells.Item(50, 6).Value = 1093
